# Weekly update: insert two new rows of fresh data at the top of the
# data block (rows 62-63), pushing the existing historical rows down by
# two positions (old row 62 -> new row 64, ... old row 89 -> new row 91).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 62; this shifts rows
# 62..89 down to 64..91 and keeps all of their cell values/styles intact.
$ws.Rows("62:63").Insert()

# Populate the two newly inserted rows with the new weekly records.
# Columns A,B,C,E,F,G,H,I,J,K,Q,T are identical across the whole sheet
# for this market/product combination.

# Row 62: Primera quality, new date 2022-02-17 (serial 44609)
$ws.Range("A62").Value = 11
$ws.Range("B62").Value = "Vega Monumental Concepción"
$ws.Range("C62").Value = "Bíobío"
$ws.Range("D62").Value = 44609
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = "Fruta"
$ws.Range("G62").Value = 100101
$ws.Range("H62").Value = "Berries"
$ws.Range("I62").Value = 100101001
$ws.Range("J62").Value = "Arándano (blue)"
$ws.Range("K62").Value = "Sin especificar"
$ws.Range("L62").Value = "Primera"
$ws.Range("M62").Value = 100
$ws.Range("N62").Value = 3500
$ws.Range("O62").Value = 4000
$ws.Range("P62").Value = 3750
$ws.Range("Q62").Value = "$/bandeja 2 kilos"
$ws.Range("R62").Value = "Región de Ñuble"
$ws.Range("S62").Value = 1875
$ws.Range("T62").Value = 2

# Row 63: Segunda quality, same new date 2022-02-17 (serial 44609)
$ws.Range("A63").Value = 11
$ws.Range("B63").Value = "Vega Monumental Concepción"
$ws.Range("C63").Value = "Bíobío"
$ws.Range("D63").Value = 44609
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = "Fruta"
$ws.Range("G63").Value = 100101
$ws.Range("H63").Value = "Berries"
$ws.Range("I63").Value = 100101001
$ws.Range("J63").Value = "Arándano (blue)"
$ws.Range("K63").Value = "Sin especificar"
$ws.Range("L63").Value = "Segunda"
$ws.Range("M63").Value = 50
$ws.Range("N63").Value = 3000
$ws.Range("O63").Value = 3000
$ws.Range("P63").Value = 3000
$ws.Range("Q63").Value = "$/bandeja 2 kilos"
$ws.Range("R63").Value = "Región de Ñuble"
$ws.Range("S63").Value = 1500
$ws.Range("T63").Value = 2
